$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet right after "总计" / right before
#    "2022-Q1". We build it by duplicating the existing "2022-Q1" sheet
#    (so it inherits identical sheet-level formatting: borders, column
#    header styles, page margins, etc.) and then overwrite its data with
#    the new quarter's figures. The original "2022-Q1" sheet is left
#    completely untouched.
# ---------------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item(2)
$wsQ1.Copy($wsQ1)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Restore the active-tab flag to the last sheet (matches the original
# workbook, where only the trailing "2020-Q4" sheet carried tabSelected).
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()

# ---------------------------------------------------------------------------
# 2) Overwrite the new "2022-Q3" sheet with the new fund-holding figures.
# ---------------------------------------------------------------------------
$wsQ3.Cells.Item(2,1).Value = 0
$wsQ3.Cells.Item(2,2).NumberFormat = "@"
$wsQ3.Cells.Item(2,2).Value = "001735"
$wsQ3.Cells.Item(2,3).Value = "广发百发大数据策略成长灵活配置混合E"
$wsQ3.Cells.Item(2,4).NumberFormat = "@"
$wsQ3.Cells.Item(2,4).Value = "0.51"
$wsQ3.Cells.Item(2,5).NumberFormat = "@"
$wsQ3.Cells.Item(2,5).Value = "43.77"
$wsQ3.Cells.Item(2,6).NumberFormat = "@"
$wsQ3.Cells.Item(2,6).Value = "2.00"
$wsQ3.Cells.Item(2,7).NumberFormat = "@"
$wsQ3.Cells.Item(2,7).Value = "0.0102"
$wsQ3.Cells.Item(2,8).Value = 2

$wsQ3.Cells.Item(3,1).Value = 1
$wsQ3.Cells.Item(3,2).NumberFormat = "@"
$wsQ3.Cells.Item(3,2).Value = "001734"
$wsQ3.Cells.Item(3,3).Value = "广发百发大数据策略成长灵活配置混合A"
$wsQ3.Cells.Item(3,4).NumberFormat = "@"
$wsQ3.Cells.Item(3,4).Value = "0.27"
$wsQ3.Cells.Item(3,5).NumberFormat = "@"
$wsQ3.Cells.Item(3,5).Value = "43.77"
$wsQ3.Cells.Item(3,6).NumberFormat = "@"
$wsQ3.Cells.Item(3,6).Value = "2.00"
$wsQ3.Cells.Item(3,7).NumberFormat = "@"
$wsQ3.Cells.Item(3,7).Value = "0.0054"
$wsQ3.Cells.Item(3,8).Value = 2

# ---------------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: insert the 2022-Q3 row at the top of
#    the data and push the existing rows down by one.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Copy format of the last existing data row down into the new row 5 first,
# so the new row's "A" cell picks up the same style as the others.
$wsTotal.Cells.Item(4,1).Copy($wsTotal.Cells.Item(5,1))

$wsTotal.Cells.Item(5,1).Value = 3
$wsTotal.Cells.Item(5,2).Value = "2020-Q4"
$wsTotal.Cells.Item(5,3).Value = 2
$wsTotal.Cells.Item(5,4).Value = 0.07000000000000001

$wsTotal.Cells.Item(4,2).Value = "2021-Q2"
$wsTotal.Cells.Item(4,3).Value = 1
$wsTotal.Cells.Item(4,4).Value = 0.03

$wsTotal.Cells.Item(3,2).Value = "2022-Q1"
$wsTotal.Cells.Item(3,3).Value = 2
$wsTotal.Cells.Item(3,4).Value = 0

$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,4).Value = 0.02
